# Update 2025-09-12: add new client "SALAZAR VERA ENRIQUE WILLIAM" into the
# alphabetically sorted client lists of the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets (inserted right before "SOLIS OCAMPO DIMAS ABDON",
# within the "OFICINA-CATAECSA" group), with all-zero sales figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": columns A:R, new row inserted at row 294
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$newRow1 = 294
$ws1.Rows.Item($newRow1).Insert()

$ws1.Cells.Item($newRow1, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item($newRow1, 2).Value = "SALAZAR VERA ENRIQUE WILLIAM"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item($newRow1, $col).Value = 0
}

# Update the trailing "N de 321" counters (now on row 324) to "N de 322"
$totalsRow1 = 324
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item($totalsRow1, $col)
    $cell.Value = $cell.Value2.ToString().Replace("de 321", "de 322")
}

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL": columns A:G, new row inserted at row 298
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$newRow2 = 298
$ws2.Rows.Item($newRow2).Insert()

$ws2.Cells.Item($newRow2, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item($newRow2, 2).Value = "SALAZAR VERA ENRIQUE WILLIAM"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item($newRow2, $col).Value = 0
}
